$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicated wording in the waste-name shared string used across the sheet
# ("Paper and cardboard packaging and cardboard packaging" -> "Paper and cardboard packaging")
$ws.Cells.Replace("Paper and cardboard packaging and cardboard packaging", "Paper and cardboard packaging")

# Leave the cursor on the cell that was corrected
$ws.Range("E3").Select()
